$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the affected ranges so numeric-looking strings
# (e.g. "328.27", "12.00", "0.110") are preserved verbatim as text,
# matching the original inlineStr cell type instead of being coerced
# into floating point numbers by COM auto-detection.
$numRng = $ws.Range("D2:E51")
$txtRng = $ws.Range("B13:C14")
$numRng.NumberFormat = "@"
$txtRng.NumberFormat = "@"

$ws.Range('D2').Value = '43.490.21'
$ws.Range('E2').Value = '  +2.80%  '
$ws.Range('D3').Value = '2.409.30'
$ws.Range('E3').Value = '  +8.48%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '328.27'
$ws.Range('E5').Value = '  +13.17%  '
$ws.Range('D6').Value = '105.12'
$ws.Range('E6').Value = '  -4.93%  '
$ws.Range('D7').Value = '0.657'
$ws.Range('E7').Value = '  +5.27%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.660'
$ws.Range('E9').Value = '  +10.57%  '
$ws.Range('D10').Value = '42.30'
$ws.Range('E10').Value = '  -2.94%  '
$ws.Range('D11').Value = '0.0946'
$ws.Range('E11').Value = '  +4.22%  '
$ws.Range('D12').Value = '8.63'
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '17.36'
$ws.Range('E13').Value = '  +16.94%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '1.03'
$ws.Range('E14').Value = '  +2.12%  '
$ws.Range('E15').Value = '  +3.16%  '
$ws.Range('D16').Value = '2.774.66'
$ws.Range('E16').Value = '  +8.46%  '
$ws.Range('D17').Value = '2.418.28'
$ws.Range('E17').Value = '  +9.25%  '
$ws.Range('D18').Value = '43.484.64'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '7.49'
$ws.Range('E19').Value = '  +5.51%  '
$ws.Range('D20').Value = '0.0000109'
$ws.Range('E20').Value = '  +4.71%  '
$ws.Range('D21').Value = '76.15'
$ws.Range('E21').Value = '  +4.95%  '
$ws.Range('D22').Value = '3.53'
$ws.Range('E22').Value = '  +5.69%  '
$ws.Range('D23').Value = '271.18'
$ws.Range('E23').Value = '  +16.26%  '
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').Value = '9.71'
$ws.Range('E25').Value = '  +8.75%  '
$ws.Range('D26').Value = '12.00'
$ws.Range('E26').Value = '  +5.60%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = '22.99'
$ws.Range('E29').Value = '  +10.51%  '
$ws.Range('D30').Value = '177.81'
$ws.Range('E30').Value = '  +2.79%  '
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('D32').Value = '38.04'
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('E33').Value = '  +4.95%  '
$ws.Range('D34').Value = '0.0938'
$ws.Range('E34').Value = '  +7.02%  '
$ws.Range('E35').Value = '  +6.56%  '
$ws.Range('E36').Value = '  +7.11%  '
$ws.Range('D37').Value = '4.90'
$ws.Range('E37').Value = '  -0.82%  '
$ws.Range('D38').Value = '4.11'
$ws.Range('E38').Value = '  -2.47%  '
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').Value = '0.110'
$ws.Range('E40').Value = '  +5.76%  '
$ws.Range('D41').Value = '2.90'
$ws.Range('E41').Value = '  +21.62%  '
$ws.Range('E42').Value = '  +24.22%  '
$ws.Range('D43').Value = '127.51'
$ws.Range('E43').Value = '  +26.56%  '
$ws.Range('D44').Value = '0.235'
$ws.Range('E44').Value = '  +2.53%  '
$ws.Range('D45').Value = '70.17'
$ws.Range('E45').Value = '  -3.34%  '
$ws.Range('E46').Value = '  +4.39%  '
$ws.Range('D48').Value = '9.78'
$ws.Range('E48').Value = '  +16.64%  '
$ws.Range('E49').Value = '  +7.30%  '
$ws.Range('D50').Value = '88.57'
$ws.Range('E50').Value = '  +65.39%  '
$ws.Range('D51').Value = '1.33'
$ws.Range('E51').Value = '  +4.54%  '

# Restore the cells to their original (default) style now that the
# text values are committed, so no stray style index is left behind.
$numRng.ClearFormats()
$txtRng.ClearFormats()
